$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2013201320132013
$ws.Range("C2").Value = 0.5610561056105611
$ws.Range("J2").Value = 0.0033003300330033
$ws.Range("P2").Value = 0.1452145214521452
$ws.Range("S2").Value = 0.0891089108910891
$ws.Range("B3").Value = 0.005813953488372093
$ws.Range("C3").Value = 0.02906976744186046
$ws.Range("J3").Value = 0.01744186046511628
$ws.Range("P3").Value = 0.7151162790697675
$ws.Range("S3").Value = 0.2325581395348837
$ws.Range("J4").Value = 0.1219512195121951
$ws.Range("P4").Value = 0.5365853658536586
$ws.Range("S4").Value = 0.3414634146341464
$ws.Range("B6").Value = 0.04830917874396135
$ws.Range("D6").Value = 0.004830917874396135
$ws.Range("F6").Value = 0.05314009661835749
$ws.Range("J6").Value = 0.178743961352657
$ws.Range("O6").Value = 0.02898550724637681
$ws.Range("Q6").Value = 0.1449275362318841
$ws.Range("R6").Value = 0.08695652173913043
$ws.Range("S6").Value = 0.4541062801932367
$ws.Range("B7").Value = 0.1483870967741935
$ws.Range("D7").Value = 0.01290322580645161
$ws.Range("F7").Value = 0.06451612903225806
$ws.Range("J7").Value = 0.08387096774193549
$ws.Range("O7").Value = 0.01290322580645161
$ws.Range("Q7").Value = 0.1741935483870968
$ws.Range("R7").Value = 0.06451612903225806
$ws.Range("S7").Value = 0.4387096774193548
$ws.Range("B8").Value = 0.115839243498818
$ws.Range("D8").Value = 0.01418439716312057
$ws.Range("F8").Value = 0.05437352245862884
$ws.Range("J8").Value = 0.1016548463356974
$ws.Range("O8").Value = 0.004728132387706856
$ws.Range("Q8").Value = 0.1631205673758865
$ws.Range("R8").Value = 0.09456264775413711
$ws.Range("S8").Value = 0.4515366430260047
$ws.Range("B9").Value = 0.09333333333333334
$ws.Range("D9").Value = 0.01777777777777778
$ws.Range("F9").Value = 0.08888888888888889
$ws.Range("J9").Value = 0.1377777777777778
$ws.Range("O9").Value = 0.01333333333333333
$ws.Range("Q9").Value = 0.1955555555555556
$ws.Range("R9").Value = 0.07555555555555556
$ws.Range("S9").Value = 0.3777777777777778
$ws.Range("B10").Value = 0.1104699093157461
$ws.Range("D10").Value = 0.02390766694146744
$ws.Range("E10").Value = 0.003297609233305853
$ws.Range("F10").Value = 0.07502061005770816
$ws.Range("J10").Value = 0.1112943116240726
$ws.Range("O10").Value = 0.01071723000824402
$ws.Range("Q10").Value = 0.2019785655399835
$ws.Range("R10").Value = 0.0898598516075845
$ws.Range("S10").Value = 0.3734542456718879
$ws.Range("G11").Value = 0.1764705882352941
$ws.Range("J11").Value = 0.09243697478991597
$ws.Range("K11").Value = 0.226890756302521
$ws.Range("L11").Value = 0.4957983193277311
$ws.Range("S11").Value = 0.008403361344537815
$ws.Range("G12").Value = 0.7235772357723578
$ws.Range("J12").Value = 0.1788617886178862
$ws.Range("K12").Value = 0.03252032520325204
$ws.Range("L12").Value = 0.03252032520325204
$ws.Range("S12").Value = 0.03252032520325204
$ws.Range("G13").Value = 0.6170212765957447
$ws.Range("J13").Value = 0.3191489361702128
$ws.Range("S13").Value = 0.06382978723404255
$ws.Range("H15").Value = 0.1707317073170732
$ws.Range("I15").Value = 0.07804878048780488
$ws.Range("J15").Value = 0.4195121951219512
$ws.Range("K15").Value = 0.08292682926829269
$ws.Range("M15").Value = 0.00975609756097561
$ws.Range("O15").Value = 0.04878048780487805
$ws.Range("S15").Value = 0.1902439024390244
$ws.Range("F16").Value = 0.005524861878453038
$ws.Range("H16").Value = 0.143646408839779
$ws.Range("I16").Value = 0.1160220994475138
$ws.Range("J16").Value = 0.4530386740331492
$ws.Range("K16").Value = 0.05524861878453038
$ws.Range("M16").Value = 0.02762430939226519
$ws.Range("O16").Value = 0.04972375690607735
$ws.Range("S16").Value = 0.1491712707182321
$ws.Range("F17").Value = 0.02403846153846154
$ws.Range("H17").Value = 0.1682692307692308
$ws.Range("I17").Value = 0.09375
$ws.Range("J17").Value = 0.4278846153846154
$ws.Range("K17").Value = 0.06009615384615385
$ws.Range("M17").Value = 0.01682692307692308
$ws.Range("O17").Value = 0.05288461538461538
$ws.Range("S17").Value = 0.15625
$ws.Range("F18").Value = 0.02061855670103093
$ws.Range("H18").Value = 0.2010309278350516
$ws.Range("I18").Value = 0.09793814432989691
$ws.Range("J18").Value = 0.4020618556701031
$ws.Range("K18").Value = 0.08247422680412371
$ws.Range("M18").Value = 0.0154639175257732
$ws.Range("O18").Value = 0.03608247422680412
$ws.Range("S18").Value = 0.1443298969072165
$ws.Range("F19").Value = 0.006216006216006216
$ws.Range("H19").Value = 0.1965811965811966
$ws.Range("I19").Value = 0.1002331002331002
$ws.Range("J19").Value = 0.372960372960373
$ws.Range("K19").Value = 0.08547008547008547
$ws.Range("M19").Value = 0.02408702408702409
$ws.Range("N19").Value = 0.000777000777000777
$ws.Range("O19").Value = 0.08702408702408702
$ws.Range("S19").Value = 0.1266511266511267
